# Applies the "Saldo" export update:
#   1. Remove the CONTEL (004589191) balance row.
#   2. Insert a new VALERIA (005440756) balance row just above the DAVID (004475395) row.
#   3. Replace the SYLVERSON (005683532) row with the ALPHASITIO (005305448) row, in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the CONTEL row (account 004589191) ---
$contelRow = $ws.Columns(1).Find("004589191").Row
$ws.Rows($contelRow).EntireRow.Delete() | Out-Null

# --- 2) Insert a new row for VALERIA right above DAVID's row (004475395) ---
$davidRow = $ws.Columns(1).Find("004475395").Row
$ws.Rows($davidRow).EntireRow.Insert() | Out-Null

$newRow = $davidRow
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "005440756"
$ws.Range("B$newRow").Value = "VALERIA"
$ws.Range("C$newRow").Value = 8066.73

# --- 3) Replace the SYLVERSON row (account 005683532) with ALPHASITIO's data, in place ---
$sylversonRow = $ws.Columns(1).Find("005683532").Row
$ws.Range("A$sylversonRow").NumberFormat = "@"
$ws.Range("A$sylversonRow").Value = "005305448"
$ws.Range("B$sylversonRow").Value = "ALPHASITIO"
$ws.Range("C$sylversonRow").Value = 1201.26
